$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-27 Monday", 2) | Out-Null
$d.Content.Find.Execute("48-24=24", $true, $false, $false, $false, $false, $true, 1, $false, "2+61=63", 2) | Out-Null
$d.Content.Find.Execute("8+38=46", $true, $false, $false, $false, $false, $true, 1, $false, "55-39=16", 2) | Out-Null
$d.Content.Find.Execute("82-13=69", $true, $false, $false, $false, $false, $true, 1, $false, "97-61=36", 2) | Out-Null
$d.Content.Find.Execute("75-62=13", $true, $false, $false, $false, $false, $true, 1, $false, "2+94=96", 2) | Out-Null
$d.Content.Find.Execute("26-7=19", $true, $false, $false, $false, $false, $true, 1, $false, "28+50=78", 2) | Out-Null
$d.Content.Find.Execute("27+48=75", $true, $false, $false, $false, $false, $true, 1, $false, "29+25=54", 2) | Out-Null
$d.Content.Find.Execute("61+33=94", $true, $false, $false, $false, $false, $true, 1, $false, "37+48=85", 2) | Out-Null
$d.Content.Find.Execute("10+57=67", $true, $false, $false, $false, $false, $true, 1, $false, "26+64=90", 2) | Out-Null
$d.Content.Find.Execute("84-6=78", $true, $false, $false, $false, $false, $true, 1, $false, "92-38=54", 2) | Out-Null
$d.Content.Find.Execute("12+82=94", $true, $false, $false, $false, $false, $true, 1, $false, "98-20=78", 2) | Out-Null
$d.Content.Find.Execute("85-39=46", $true, $false, $false, $false, $false, $true, 1, $false, "84-69=15", 2) | Out-Null
$d.Content.Find.Execute("29+62=91", $true, $false, $false, $false, $false, $true, 1, $false, "63-62=1", 2) | Out-Null
$d.Content.Find.Execute("83-72=11", $true, $false, $false, $false, $false, $true, 1, $false, "97-96=1", 2) | Out-Null
$d.Content.Find.Execute("72-42=30", $true, $false, $false, $false, $false, $true, 1, $false, "12+28=40", 2) | Out-Null
$d.Content.Find.Execute("59-38=21", $true, $false, $false, $false, $false, $true, 1, $false, "26+18=44", 2) | Out-Null
$d.Content.Find.Execute("70-54=16", $true, $false, $false, $false, $false, $true, 1, $false, "23+15=38", 2) | Out-Null
$d.Content.Find.Execute("37-1=36", $true, $false, $false, $false, $false, $true, 1, $false, "19+50=69", 2) | Out-Null
$d.Content.Find.Execute("60-15=45", $true, $false, $false, $false, $false, $true, 1, $false, "47-45=2", 2) | Out-Null
$d.Content.Find.Execute("57-44=13", $true, $false, $false, $false, $false, $true, 1, $false, "15-1=14", 2) | Out-Null
$d.Content.Find.Execute("55-6=49", $true, $false, $false, $false, $false, $true, 1, $false, "54-2=52", 2) | Out-Null
$d.Content.Find.Execute("65-60=5", $true, $false, $false, $false, $false, $true, 1, $false, "20+26=46", 2) | Out-Null
$d.Content.Find.Execute("68+6=74", $true, $false, $false, $false, $false, $true, 1, $false, "40+23=63", 2) | Out-Null
$d.Content.Find.Execute("36+62=98", $true, $false, $false, $false, $false, $true, 1, $false, "51-16=35", 2) | Out-Null
$d.Content.Find.Execute("89-20=69", $true, $false, $false, $false, $false, $true, 1, $false, "18+11=29", 2) | Out-Null
$d.Content.Find.Execute("92-86=6", $true, $false, $false, $false, $false, $true, 1, $false, "46-43=3", 2) | Out-Null
$d.Content.Find.Execute("20+65=85", $true, $false, $false, $false, $false, $true, 1, $false, "33+28=61", 2) | Out-Null
$d.Content.Find.Execute("78-40=38", $true, $false, $false, $false, $false, $true, 1, $false, "35-25=10", 2) | Out-Null
$d.Content.Find.Execute("23-12=11", $true, $false, $false, $false, $false, $true, 1, $false, "81+9=90", 2) | Out-Null
$d.Content.Find.Execute("37-36=1", $true, $false, $false, $false, $false, $true, 1, $false, "18-3=15", 2) | Out-Null
$d.Content.Find.Execute("9+23=32", $true, $false, $false, $false, $false, $true, 1, $false, "30+35=65", 2) | Out-Null
$d.Content.Find.Execute("7+3=10", $true, $false, $false, $false, $false, $true, 1, $false, "93+1=94", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=63", 2) | Out-Null
$d.Content.Find.Execute("55-14=41", $true, $false, $false, $false, $false, $true, 1, $false, "51+0=51", 2) | Out-Null
$d.Content.Find.Execute("52-34=18", $true, $false, $false, $false, $false, $true, 1, $false, "61-44=17", 2) | Out-Null
$d.Content.Find.Execute("22+1=23", $true, $false, $false, $false, $false, $true, 1, $false, "32+1=33", 2) | Out-Null
$d.Content.Find.Execute("76-50=26", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=53", 2) | Out-Null
$d.Content.Find.Execute("20-19=1", $true, $false, $false, $false, $false, $true, 1, $false, "67-61=6", 2) | Out-Null
$d.Content.Find.Execute("80+0=80", $true, $false, $false, $false, $false, $true, 1, $false, "9+8=17", 2) | Out-Null
$d.Content.Find.Execute("83-48=35", $true, $false, $false, $false, $false, $true, 1, $false, "94-31=63", 2) | Out-Null
$d.Content.Find.Execute("28+9=37", $true, $false, $false, $false, $false, $true, 1, $false, "96-17=79", 2) | Out-Null
$d.Content.Find.Execute("10-4=6", $true, $false, $false, $false, $false, $true, 1, $false, "46+37=83", 2) | Out-Null
$d.Content.Find.Execute("64-37=27", $true, $false, $false, $false, $false, $true, 1, $false, "8-4=4", 2) | Out-Null
$d.Content.Find.Execute("60-18=42", $true, $false, $false, $false, $false, $true, 1, $false, "5+51=56", 2) | Out-Null
$d.Content.Find.Execute("10+79=89", $true, $false, $false, $false, $false, $true, 1, $false, "95-71=24", 2) | Out-Null
$d.Content.Find.Execute("0+31=31", $true, $false, $false, $false, $false, $true, 1, $false, "49-21=28", 2) | Out-Null
$d.Content.Find.Execute("90-69=21", $true, $false, $false, $false, $false, $true, 1, $false, "85+11=96", 2) | Out-Null
$d.Content.Find.Execute("73-64=9", $true, $false, $false, $false, $false, $true, 1, $false, "64+33=97", 2) | Out-Null
$d.Content.Find.Execute("62-22=40", $true, $false, $false, $false, $false, $true, 1, $false, "3+91=94", 2) | Out-Null
$d.Content.Find.Execute("28+8=36", $true, $false, $false, $false, $false, $true, 1, $false, "19+3=22", 2) | Out-Null
$d.Content.Find.Execute("11+68=79", $true, $false, $false, $false, $false, $true, 1, $false, "46+18=64", 2) | Out-Null
$d.Content.Find.Execute("70-50=20", $true, $false, $false, $false, $false, $true, 1, $false, "95-45=50", 2) | Out-Null
$d.Content.Find.Execute("37+3=40", $true, $false, $false, $false, $false, $true, 1, $false, "14+12=26", 2) | Out-Null
$d.Content.Find.Execute("42+19=61", $true, $false, $false, $false, $false, $true, 1, $false, "86-57=29", 2) | Out-Null
$d.Content.Find.Execute("25-11=14", $true, $false, $false, $false, $false, $true, 1, $false, "22+60=82", 2) | Out-Null
$d.Content.Find.Execute("96-31=65", $true, $false, $false, $false, $false, $true, 1, $false, "43+45=88", 2) | Out-Null
$d.Content.Find.Execute("47+9=56", $true, $false, $false, $false, $false, $true, 1, $false, "83+8=91", 2) | Out-Null
$d.Content.Find.Execute("26+6=32", $true, $false, $false, $false, $false, $true, 1, $false, "61+16=77", 2) | Out-Null
$d.Content.Find.Execute("24-20=4", $true, $false, $false, $false, $false, $true, 1, $false, "65+17=82", 2) | Out-Null
$d.Content.Find.Execute("41-27=14", $true, $false, $false, $false, $false, $true, 1, $false, "85-35=50", 2) | Out-Null
$d.Content.Find.Execute("77-4=73", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=84", 2) | Out-Null
$d.Content.Find.Execute("17+76=93", $true, $false, $false, $false, $false, $true, 1, $false, "33+2=35", 2) | Out-Null
$d.Content.Find.Execute("54+45=99", $true, $false, $false, $false, $false, $true, 1, $false, "36-0=36", 2) | Out-Null
$d.Content.Find.Execute("19+53=72", $true, $false, $false, $false, $false, $true, 1, $false, "56-24=32", 2) | Out-Null
$d.Content.Find.Execute("8+75=83", $true, $false, $false, $false, $false, $true, 1, $false, "78-9=69", 2) | Out-Null
$d.Content.Find.Execute("92-39=53", $true, $false, $false, $false, $false, $true, 1, $false, "2+80=82", 2) | Out-Null
$d.Content.Find.Execute("76-40=36", $true, $false, $false, $false, $false, $true, 1, $false, "23+42=65", 2) | Out-Null
$d.Content.Find.Execute("53-47=6", $true, $false, $false, $false, $false, $true, 1, $false, "56-16=40", 2) | Out-Null
$d.Content.Find.Execute("43+34=77", $true, $false, $false, $false, $false, $true, 1, $false, "5+60=65", 2) | Out-Null
$d.Content.Find.Execute("57+19=76", $true, $false, $false, $false, $false, $true, 1, $false, "89-12=77", 2) | Out-Null
$d.Content.Find.Execute("59-23=36", $true, $false, $false, $false, $false, $true, 1, $false, "54+15=69", 2) | Out-Null
$d.Content.Find.Execute("70-21=49", $true, $false, $false, $false, $false, $true, 1, $false, "60-27=33", 2) | Out-Null
$d.Content.Find.Execute("93-14=79", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=73", 2) | Out-Null
$d.Content.Find.Execute("30+6=36", $true, $false, $false, $false, $false, $true, 1, $false, "9+70=79", 2) | Out-Null
$d.Content.Find.Execute("10+35=45", $true, $false, $false, $false, $false, $true, 1, $false, "14+59=73", 2) | Out-Null
$d.Content.Find.Execute("26+73=99", $true, $false, $false, $false, $false, $true, 1, $false, "46+9=55", 2) | Out-Null
$d.Content.Find.Execute("18+15=33", $true, $false, $false, $false, $false, $true, 1, $false, "62+2=64", 2) | Out-Null
$d.Content.Find.Execute("66-22=44", $true, $false, $false, $false, $false, $true, 1, $false, "29+38=67", 2) | Out-Null
$d.Content.Find.Execute("16+41=57", $true, $false, $false, $false, $false, $true, 1, $false, "43-15=28", 2) | Out-Null
$d.Content.Find.Execute("8+62=70", $true, $false, $false, $false, $false, $true, 1, $false, "42+49=91", 2) | Out-Null
$d.Content.Find.Execute("6+3=9", $true, $false, $false, $false, $false, $true, 1, $false, "5+47=52", 2) | Out-Null
$d.Content.Find.Execute("27+9=36", $true, $false, $false, $false, $false, $true, 1, $false, "66+10=76", 2) | Out-Null
$d.Content.Find.Execute("27-21=6", $true, $false, $false, $false, $false, $true, 1, $false, "23+7=30", 2) | Out-Null
$d.Content.Find.Execute("55-42=13", $true, $false, $false, $false, $false, $true, 1, $false, "12+77=89", 2) | Out-Null
$d.Content.Find.Execute("14-7=7", $true, $false, $false, $false, $false, $true, 1, $false, "65-18=47", 2) | Out-Null
$d.Content.Find.Execute("94-32=62", $true, $false, $false, $false, $false, $true, 1, $false, "62-28=34", 2) | Out-Null
$d.Content.Find.Execute("91-5=86", $true, $false, $false, $false, $false, $true, 1, $false, "67-0=67", 2) | Out-Null
$d.Content.Find.Execute("6+89=95", $true, $false, $false, $false, $false, $true, 1, $false, "11+17=28", 2) | Out-Null
$d.Content.Find.Execute("80-54=26", $true, $false, $false, $false, $false, $true, 1, $false, "96-69=27", 2) | Out-Null
$d.Content.Find.Execute("20-0=20", $true, $false, $false, $false, $false, $true, 1, $false, "48-48=0", 2) | Out-Null
$d.Content.Find.Execute("93-24=69", $true, $false, $false, $false, $false, $true, 1, $false, "2+78=80", 2) | Out-Null
$d.Content.Find.Execute("64-8=56", $true, $false, $false, $false, $false, $true, 1, $false, "70-28=42", 2) | Out-Null
$d.Content.Find.Execute("0+14=14", $true, $false, $false, $false, $false, $true, 1, $false, "14+0=14", 2) | Out-Null
$d.Content.Find.Execute("36+33=69", $true, $false, $false, $false, $false, $true, 1, $false, "41+56=97", 2) | Out-Null
$d.Content.Find.Execute("34+25=59", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=45", 2) | Out-Null
$d.Content.Find.Execute("4+5=9", $true, $false, $false, $false, $false, $true, 1, $false, "92-82=10", 2) | Out-Null
$d.Content.Find.Execute("25+19=44", $true, $false, $false, $false, $false, $true, 1, $false, "75+2=77", 2) | Out-Null
$d.Content.Find.Execute("0+51=51", $true, $false, $false, $false, $false, $true, 1, $false, "91-49=42", 2) | Out-Null
$d.Content.Find.Execute("2+91=93", $true, $false, $false, $false, $false, $true, 1, $false, "21+41=62", 2) | Out-Null
$d.Content.Find.Execute("54+22=76", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=11", 2) | Out-Null
$d.Content.Find.Execute("59-4=55", $true, $false, $false, $false, $false, $true, 1, $false, "63-53=10", 2) | Out-Null
